$wb = $excel.ActiveWorkbook

# Male_50m sheet (index 2)
$ws2 = $wb.Worksheets.Item("Male_50m")

# Row 10: rename "Solum Ole Peder Uthus" -> "Ole Peder Uthus Solum"
$ws2.Range("A10").Value = "Ole Peder Uthus Solum"

# Row 11: replace old "Ole Peder Uthus Solum" duplicate entry with Tor Arne Hegvik data
$ws2.Range("A11").Value = "Tor Arne Hegvik"
$ws2.Range("B11").Value = "2.06,70"
$ws2.Range("C11").Value = 521
$ws2.Range("D11").Value = "28.04.2007"
$ws2.Range("E11").Value = "Namsos"

# Female_50m sheet (index 4)
$ws4 = $wb.Worksheets.Item("Female_50m")

# Row 11: replace "Frøydis Vatn Andersen" with "Sigrid Eldholm" and updated data
$ws4.Range("A11").Value = "Sigrid Eldholm"
$ws4.Range("B11").Value = "2.21,22"
$ws4.Range("D11").NumberFormat = "@"
$ws4.Range("D11").Value = "02.07.2022"
$ws4.Range("E11").Value = "Stockholm"
